# Bugfixed the naive forecaster component module
# The forecast vector table is shifted up by one period (the oldest
# row is dropped) and the y_0_forecast / y_1_forecast columns are
# recomputed with corrected values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the last (now-stale) observation row; this also shrinks the
# sheet dimension from A1:E53 down to A1:E52.
$ws.Rows(53).Delete()

# Rewrite the data rows (row 1 is the header) with the corrected
# dates/years (columns A, B, D) and recomputed forecast values
# (columns C, E). Cells that should be blank are cleared.
$ws.Cells.Item(2,1).Value2 = 39583; $ws.Cells.Item(2,2).Value2 = 2008; $ws.Cells.Item(2,3).Value2 = $null; $ws.Cells.Item(2,4).Value2 = 2009; $ws.Cells.Item(2,5).Value2 = $null
$ws.Cells.Item(3,1).Value2 = 39765; $ws.Cells.Item(3,2).Value2 = 2008; $ws.Cells.Item(3,3).Value2 = $null; $ws.Cells.Item(3,4).Value2 = 2009; $ws.Cells.Item(3,5).Value2 = $null
$ws.Cells.Item(4,1).Value2 = 39948; $ws.Cells.Item(4,2).Value2 = 2009; $ws.Cells.Item(4,3).Value2 = $null; $ws.Cells.Item(4,4).Value2 = 2010; $ws.Cells.Item(4,5).Value2 = $null
$ws.Cells.Item(5,1).Value2 = 40130; $ws.Cells.Item(5,2).Value2 = 2009; $ws.Cells.Item(5,3).Value2 = -3.872359107260159; $ws.Cells.Item(5,4).Value2 = 2010; $ws.Cells.Item(5,5).Value2 = $null
$ws.Cells.Item(6,1).Value2 = 40310; $ws.Cells.Item(6,2).Value2 = 2010; $ws.Cells.Item(6,3).Value2 = $null; $ws.Cells.Item(6,4).Value2 = 2011; $ws.Cells.Item(6,5).Value2 = $null
$ws.Cells.Item(7,1).Value2 = 40494; $ws.Cells.Item(7,2).Value2 = 2010; $ws.Cells.Item(7,3).Value2 = 4.530477057343663; $ws.Cells.Item(7,4).Value2 = 2011; $ws.Cells.Item(7,5).Value2 = $null
$ws.Cells.Item(8,1).Value2 = 40676; $ws.Cells.Item(8,2).Value2 = 2011; $ws.Cells.Item(8,3).Value2 = $null; $ws.Cells.Item(8,4).Value2 = 2012; $ws.Cells.Item(8,5).Value2 = $null
$ws.Cells.Item(9,1).Value2 = 40862; $ws.Cells.Item(9,2).Value2 = 2011; $ws.Cells.Item(9,3).Value2 = 6.833902841285977; $ws.Cells.Item(9,4).Value2 = 2012; $ws.Cells.Item(9,5).Value2 = $null
$ws.Cells.Item(10,1).Value2 = 41044; $ws.Cells.Item(10,2).Value2 = 2012; $ws.Cells.Item(10,3).Value2 = $null; $ws.Cells.Item(10,4).Value2 = 2013; $ws.Cells.Item(10,5).Value2 = $null
$ws.Cells.Item(11,1).Value2 = 41228; $ws.Cells.Item(11,2).Value2 = 2012; $ws.Cells.Item(11,3).Value2 = 4.166536506645224; $ws.Cells.Item(11,4).Value2 = 2013; $ws.Cells.Item(11,5).Value2 = 2.693188401769642
$ws.Cells.Item(12,1).Value2 = 41409; $ws.Cells.Item(12,2).Value2 = 2013; $ws.Cells.Item(12,3).Value2 = 1.785377844167058; $ws.Cells.Item(12,4).Value2 = 2014; $ws.Cells.Item(12,5).Value2 = 2.333075171696652
$ws.Cells.Item(13,1).Value2 = 41592; $ws.Cells.Item(13,2).Value2 = 2013; $ws.Cells.Item(13,3).Value2 = 2.669880057548091; $ws.Cells.Item(13,4).Value2 = 2014; $ws.Cells.Item(13,5).Value2 = 3.947916604971446
$ws.Cells.Item(14,1).Value2 = 41774; $ws.Cells.Item(14,2).Value2 = 2014; $ws.Cells.Item(14,3).Value2 = 5.477304442308206; $ws.Cells.Item(14,4).Value2 = 2015; $ws.Cells.Item(14,5).Value2 = 4.052456259163839
$ws.Cells.Item(15,1).Value2 = 41957; $ws.Cells.Item(15,2).Value2 = 2014; $ws.Cells.Item(15,3).Value2 = 5.50293301232252; $ws.Cells.Item(15,4).Value2 = 2015; $ws.Cells.Item(15,5).Value2 = 4.998814576944932
$ws.Cells.Item(16,1).Value2 = 42137; $ws.Cells.Item(16,2).Value2 = 2015; $ws.Cells.Item(16,3).Value2 = 4.666532690711245; $ws.Cells.Item(16,4).Value2 = 2016; $ws.Cells.Item(16,5).Value2 = 3.659383764712709
$ws.Cells.Item(17,1).Value2 = 42321; $ws.Cells.Item(17,2).Value2 = 2015; $ws.Cells.Item(17,3).Value2 = 4.829481320500406; $ws.Cells.Item(17,4).Value2 = 2016; $ws.Cells.Item(17,5).Value2 = 4.673582741620552
$ws.Cells.Item(18,1).Value2 = 42503; $ws.Cells.Item(18,2).Value2 = 2016; $ws.Cells.Item(18,3).Value2 = 5.266214435142658; $ws.Cells.Item(18,4).Value2 = 2017; $ws.Cells.Item(18,5).Value2 = 4.181342739750682
$ws.Cells.Item(19,1).Value2 = 42689; $ws.Cells.Item(19,2).Value2 = 2016; $ws.Cells.Item(19,3).Value2 = 5.100281927437122; $ws.Cells.Item(19,4).Value2 = 2017; $ws.Cells.Item(19,5).Value2 = 4.372458986620376
$ws.Cells.Item(20,1).Value2 = 42867; $ws.Cells.Item(20,2).Value2 = 2017; $ws.Cells.Item(20,3).Value2 = 4.811826107786477; $ws.Cells.Item(20,4).Value2 = 2018; $ws.Cells.Item(20,5).Value2 = 4.131858242365549
$ws.Cells.Item(21,1).Value2 = 43053; $ws.Cells.Item(21,2).Value2 = 2017; $ws.Cells.Item(21,3).Value2 = 5.161358932333737; $ws.Cells.Item(21,4).Value2 = 2018; $ws.Cells.Item(21,5).Value2 = 4.927320050172312
$ws.Cells.Item(22,1).Value2 = 43145; $ws.Cells.Item(22,2).Value2 = 2018; $ws.Cells.Item(22,3).Value2 = 6.022380124455107; $ws.Cells.Item(22,4).Value2 = 2019; $ws.Cells.Item(22,5).Value2 = 4.474956658559948
$ws.Cells.Item(23,1).Value2 = 43235; $ws.Cells.Item(23,2).Value2 = 2018; $ws.Cells.Item(23,3).Value2 = 5.91185619417105; $ws.Cells.Item(23,4).Value2 = 2019; $ws.Cells.Item(23,5).Value2 = 4.365509285986957
$ws.Cells.Item(24,1).Value2 = 43326; $ws.Cells.Item(24,2).Value2 = 2018; $ws.Cells.Item(24,3).Value2 = 5.904095356703798; $ws.Cells.Item(24,4).Value2 = 2019; $ws.Cells.Item(24,5).Value2 = 4.348199743880454
$ws.Cells.Item(25,1).Value2 = 43418; $ws.Cells.Item(25,2).Value2 = 2018; $ws.Cells.Item(25,3).Value2 = 5.902681694119694; $ws.Cells.Item(25,4).Value2 = 2019; $ws.Cells.Item(25,5).Value2 = 4.339089271348406
$ws.Cells.Item(26,1).Value2 = 43510; $ws.Cells.Item(26,2).Value2 = 2019; $ws.Cells.Item(26,3).Value2 = 3.109393707322261; $ws.Cells.Item(26,4).Value2 = 2020; $ws.Cells.Item(26,5).Value2 = 3.558392386986431
$ws.Cells.Item(27,1).Value2 = 43600; $ws.Cells.Item(27,2).Value2 = 2019; $ws.Cells.Item(27,3).Value2 = 5.114185474093769; $ws.Cells.Item(27,4).Value2 = 2020; $ws.Cells.Item(27,5).Value2 = 5.472991335528654
$ws.Cells.Item(28,1).Value2 = 43691; $ws.Cells.Item(28,2).Value2 = 2019; $ws.Cells.Item(28,3).Value2 = 3.799522169175473; $ws.Cells.Item(28,4).Value2 = 2020; $ws.Cells.Item(28,5).Value2 = 2.778402897289434
$ws.Cells.Item(29,1).Value2 = 43783; $ws.Cells.Item(29,2).Value2 = 2019; $ws.Cells.Item(29,3).Value2 = 3.884502719230132; $ws.Cells.Item(29,4).Value2 = 2020; $ws.Cells.Item(29,5).Value2 = 3.243024666552685
$ws.Cells.Item(30,1).Value2 = 43875; $ws.Cells.Item(30,2).Value2 = 2020; $ws.Cells.Item(30,3).Value2 = 3.662599762249985; $ws.Cells.Item(30,4).Value2 = 2021; $ws.Cells.Item(30,5).Value2 = 3.993000457359908
$ws.Cells.Item(31,1).Value2 = 43966; $ws.Cells.Item(31,2).Value2 = 2020; $ws.Cells.Item(31,3).Value2 = 2.167530781895133; $ws.Cells.Item(31,4).Value2 = 2021; $ws.Cells.Item(31,5).Value2 = 2.573593955528963
$ws.Cells.Item(32,1).Value2 = 44068; $ws.Cells.Item(32,2).Value2 = 2020; $ws.Cells.Item(32,3).Value2 = -3.840397826549158; $ws.Cells.Item(32,4).Value2 = 2021; $ws.Cells.Item(32,5).Value2 = -5.097705497973837
$ws.Cells.Item(33,1).Value2 = 44159; $ws.Cells.Item(33,2).Value2 = 2020; $ws.Cells.Item(33,3).Value2 = -3.840397826549158; $ws.Cells.Item(33,4).Value2 = 2021; $ws.Cells.Item(33,5).Value2 = 0.2915162802050064
$ws.Cells.Item(34,1).Value2 = 44251; $ws.Cells.Item(34,2).Value2 = 2021; $ws.Cells.Item(34,3).Value2 = -0.1964516829170981; $ws.Cells.Item(34,4).Value2 = 2022; $ws.Cells.Item(34,5).Value2 = 3.329288211255621
$ws.Cells.Item(35,1).Value2 = 44341; $ws.Cells.Item(35,2).Value2 = 2021; $ws.Cells.Item(35,3).Value2 = 0.5766229317536675; $ws.Cells.Item(35,4).Value2 = 2022; $ws.Cells.Item(35,5).Value2 = 4.059584075094214
$ws.Cells.Item(36,1).Value2 = 44432; $ws.Cells.Item(36,2).Value2 = 2021; $ws.Cells.Item(36,3).Value2 = 0.4839811651348835; $ws.Cells.Item(36,4).Value2 = 2022; $ws.Cells.Item(36,5).Value2 = 3.96063514023246
$ws.Cells.Item(37,1).Value2 = 44525; $ws.Cells.Item(37,2).Value2 = 2021; $ws.Cells.Item(37,3).Value2 = 0.4839811651348835; $ws.Cells.Item(37,4).Value2 = 2022; $ws.Cells.Item(37,5).Value2 = 3.818597641626909
$ws.Cells.Item(38,1).Value2 = 44617; $ws.Cells.Item(38,2).Value2 = 2022; $ws.Cells.Item(38,3).Value2 = 2.256289400228262; $ws.Cells.Item(38,4).Value2 = 2023; $ws.Cells.Item(38,5).Value2 = 3.516729866534796
$ws.Cells.Item(39,1).Value2 = 44706; $ws.Cells.Item(39,2).Value2 = 2022; $ws.Cells.Item(39,3).Value2 = 2.288114387968587; $ws.Cells.Item(39,4).Value2 = 2023; $ws.Cells.Item(39,5).Value2 = 3.463553906111505
$ws.Cells.Item(40,1).Value2 = 44798; $ws.Cells.Item(40,2).Value2 = 2022; $ws.Cells.Item(40,3).Value2 = 2.06342951900429; $ws.Cells.Item(40,4).Value2 = 2023; $ws.Cells.Item(40,5).Value2 = 3.020376488332777
$ws.Cells.Item(41,1).Value2 = 44890; $ws.Cells.Item(41,2).Value2 = 2022; $ws.Cells.Item(41,3).Value2 = 2.06342951900429; $ws.Cells.Item(41,4).Value2 = 2023; $ws.Cells.Item(41,5).Value2 = 0.6985632195332103
$ws.Cells.Item(42,1).Value2 = 44981; $ws.Cells.Item(42,2).Value2 = 2023; $ws.Cells.Item(42,3).Value2 = -1.252226393550548; $ws.Cells.Item(42,4).Value2 = 2024; $ws.Cells.Item(42,5).Value2 = 2.967032781824974
$ws.Cells.Item(43,1).Value2 = 45071; $ws.Cells.Item(43,2).Value2 = 2023; $ws.Cells.Item(43,3).Value2 = -2.013802094285932; $ws.Cells.Item(43,4).Value2 = 2024; $ws.Cells.Item(43,5).Value2 = 2.374210810973465
$ws.Cells.Item(44,1).Value2 = 45163; $ws.Cells.Item(44,2).Value2 = 2023; $ws.Cells.Item(44,3).Value2 = -2.156362896191677; $ws.Cells.Item(44,4).Value2 = 2024; $ws.Cells.Item(44,5).Value2 = 2.034789645219792
$ws.Cells.Item(45,1).Value2 = 45254; $ws.Cells.Item(45,2).Value2 = 2023; $ws.Cells.Item(45,3).Value2 = -2.156362896191677; $ws.Cells.Item(45,4).Value2 = 2024; $ws.Cells.Item(45,5).Value2 = 0.3452735157291054
$ws.Cells.Item(46,1).Value2 = 45345; $ws.Cells.Item(46,2).Value2 = 2024; $ws.Cells.Item(46,3).Value2 = -0.4399034310282546; $ws.Cells.Item(46,4).Value2 = 2025; $ws.Cells.Item(46,5).Value2 = 2.563033601911258
$ws.Cells.Item(47,1).Value2 = 45436; $ws.Cells.Item(47,2).Value2 = 2024; $ws.Cells.Item(47,3).Value2 = -0.5865622195987186; $ws.Cells.Item(47,4).Value2 = 2025; $ws.Cells.Item(47,5).Value2 = 2.431929210693595
$ws.Cells.Item(48,1).Value2 = 45534; $ws.Cells.Item(48,2).Value2 = 2024; $ws.Cells.Item(48,3).Value2 = -0.8205034771073372; $ws.Cells.Item(48,4).Value2 = 2025; $ws.Cells.Item(48,5).Value2 = 1.969879323458756
$ws.Cells.Item(49,1).Value2 = 45618; $ws.Cells.Item(49,2).Value2 = 2024; $ws.Cells.Item(49,3).Value2 = -0.8205034771073372; $ws.Cells.Item(49,4).Value2 = 2025; $ws.Cells.Item(49,5).Value2 = 1.5902148106679
$ws.Cells.Item(50,1).Value2 = 45713; $ws.Cells.Item(50,2).Value2 = 2025; $ws.Cells.Item(50,3).Value2 = 1.06642809951869; $ws.Cells.Item(50,4).Value2 = 2026; $ws.Cells.Item(50,5).Value2 = 2.36642828939615
$ws.Cells.Item(51,1).Value2 = 45800; $ws.Cells.Item(51,2).Value2 = 2025; $ws.Cells.Item(51,3).Value2 = 0.7174582534189566; $ws.Cells.Item(51,4).Value2 = 2026; $ws.Cells.Item(51,5).Value2 = 2.061048937680932
$ws.Cells.Item(52,1).Value2 = 45891; $ws.Cells.Item(52,2).Value2 = 2025; $ws.Cells.Item(52,3).Value2 = 0.5750555200350504; $ws.Cells.Item(52,4).Value2 = 2026; $ws.Cells.Item(52,5).Value2 = 1.743520202089877
